# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (with per-fund holdings) right after
#   the existing "2021-Q4" sheet and before "总计".
# - Update the "总计" (summary) sheet with a new first data row for
#   2022-Q1, pushing the existing 2021-Q4 summary row down.
#
# NOTE: worksheet object references captured before a
# `Worksheets.Add()` call can go stale (they re-resolve by their
# original numeric position once the sheet collection changes), so every
# sheet handle below is (re)fetched by name immediately before use.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($cell, [string]$text)
    # Force literal text storage (avoids Excel re-interpreting numeric-
    # looking strings like "012751" / "0.34" as numbers), then strip the
    # number-format override so the cell keeps the workbook's default style.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right after "2021-Q4".
# ---------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "2022-Q1"

$sheetQ1 = $wb.Worksheets.Item("2022-Q1")

# Header row
$sheetQ1.Cells.Item(1, 2).Value = "基金代码"
$sheetQ1.Cells.Item(1, 3).Value = "基金名称"
$sheetQ1.Cells.Item(1, 4).Value = "基金规模"
$sheetQ1.Cells.Item(1, 5).Value = "股票总仓位"
$sheetQ1.Cells.Item(1, 6).Value = "仓位占比"
$sheetQ1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$sheetQ1.Cells.Item(1, 8).Value = "仓位排名"

$styleSrc = $wb.Worksheets.Item("2021-Q4").Cells.Item(2, 1)
$styleSrc.Copy()
$wb.Worksheets.Item("2022-Q1").Range("B1:H1").PasteSpecial(-4122)

# Data rows
$fundRows = @(
    @{ Row = 2; Idx = 0; Code = "012751"; Name = "建信纳斯达克100指数（QDII）A 美元现汇"; Size = "0.34"; Position = "88.02"; Ratio = "1.83"; Value = "0.0062"; Rank = 9 },
    @{ Row = 3; Idx = 1; Code = "012752"; Name = "建信纳斯达克100指数（QDII）C 人民币";     Size = "0.34"; Position = "88.02"; Ratio = "1.83"; Value = "0.0062"; Rank = 9 },
    @{ Row = 4; Idx = 2; Code = "012753"; Name = "建信纳斯达克100指数（QDII）C 美元现汇"; Size = "0.34"; Position = "88.02"; Ratio = "1.83"; Value = "0.0062"; Rank = 9 }
)

foreach ($r in $fundRows) {
    $sheetQ1 = $wb.Worksheets.Item("2022-Q1")
    $sheetQ1.Cells.Item($r.Row, 1).Value = $r.Idx
    Set-TextValue $sheetQ1.Cells.Item($r.Row, 2) $r.Code
    Set-TextValue $sheetQ1.Cells.Item($r.Row, 3) $r.Name
    Set-TextValue $sheetQ1.Cells.Item($r.Row, 4) $r.Size
    Set-TextValue $sheetQ1.Cells.Item($r.Row, 5) $r.Position
    Set-TextValue $sheetQ1.Cells.Item($r.Row, 6) $r.Ratio
    Set-TextValue $sheetQ1.Cells.Item($r.Row, 7) $r.Value
    $sheetQ1.Cells.Item($r.Row, 8).Value = $r.Rank
}

$styleSrc = $wb.Worksheets.Item("2021-Q4").Cells.Item(2, 1)
$styleSrc.Copy()
$wb.Worksheets.Item("2022-Q1").Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing 2021-Q4 row to row 3 and
#    insert a new 2022-Q1 row at row 2.
# ---------------------------------------------------------------------
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetTotal.Cells.Item(3, 1).Value = 1
Set-TextValue $sheetTotal.Cells.Item(3, 2) "2021-Q4"
$sheetTotal.Cells.Item(3, 3).Value = 3
$sheetTotal.Cells.Item(3, 4).Value = 0.38

$sheetTotal = $wb.Worksheets.Item("总计")
$sheetTotal.Cells.Item(2, 1).Value = 0
Set-TextValue $sheetTotal.Cells.Item(2, 2) "2022-Q1"
$sheetTotal.Cells.Item(2, 3).Value = 3
$sheetTotal.Cells.Item(2, 4).Value = 0.02

$styleSrc = $wb.Worksheets.Item("2021-Q4").Cells.Item(2, 1)
$styleSrc.Copy()
$sheetTotal = $wb.Worksheets.Item("总计")
$sheetTotal.Range("A2:A3").PasteSpecial(-4122)

Write-Host "2022-Q1 sheet inserted and summary sheet updated"
